$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Fix "ingrating" -> "integrating" in the disaster-recovery-plan paragraph.
# ---------------------------------------------------------------------------
$p39 = $d.Paragraphs(39)
$p39.Range.Find.Execute("ingrating", $false, $false, $false, $false, $false, `
    $true, 1, $false, "integrating", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Expand the "Anyone using a SQL server database..." paragraph with the
#    new "Many companies cannot function..." sentence.
# ---------------------------------------------------------------------------
$p41 = $d.Paragraphs(41)
$p41.Range.Find.Execute( `
    "after a disaster. DR plans are more than just backing up data,", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "after a disaster. Many companies cannot function without or afford to lose even a small amount of data. DR plans are more than just backing up data,", `
    2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Fill in the blank paragraph before "5. Find two companies..." with the
#    Plan B / Zetta comparison text, then add the follow-up paragraphs:
#    - a paragraph comparing the two services in more depth
#    - a blank line
#    - a hyperlink to the techradar article
#    - a blank line
#    - a hyperlink to zetta.net
#    - a blank line
#    - a hyperlink to planb.co.uk
# ---------------------------------------------------------------------------

$rsquo = [char]8217
$lsquo = [char]8216

# Locate the blank paragraph that immediately precedes "5.  Find two companies"
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "5.*Find*companies*disaster recovery services*") {
        $targetIndex = $i - 1
        break
    }
}

$pC = $d.Paragraphs($targetIndex)
$rC = $pC.Range
$rC.Collapse(1)
$rC.InsertAfter("Plan B Disaster Recovery and Zetta Disaster Recovery are providers of disaster recovery solutions for SQL servers. Arcserve provides a free trial, however Plan B does not. Both provide instant recovery and protection for both physical and virtual servers on public and private clouds. With Plan B, you have access to specialist engineers who can provide tech support and help. Plan B also tests the replica system every day for flaws. They also offer custom solution designs.")

# Paragraph D: Zetta discussion
$pC2 = $d.Paragraphs($targetIndex)
$pC2.Range.InsertParagraphAfter()
$pD = $d.Paragraphs($targetIndex + 1)
$rD = $pD.Range
$rD.Collapse(1)
$rD.InsertAfter("Zetta Disaster Recovery may be the better option for a less technical user with it" + $rsquo + "s simple and effect " + $lsquo + "push-button" + $rsquo + " recovery. This product also looks to be the more cost-effective for smaller businesses that may not be able to afford a more pricey disaster recovery service. Zetta also offers experts to assist in the initial setup and testing of the service. Their engineers are available to help during failback to make sure the transition after a disaster is as smooth as possible. Zetta uses cloud-based disaster recovery requiring no on-premise device.")

# Paragraph E: blank
$pD2 = $d.Paragraphs($targetIndex + 1)
$pD2.Range.InsertParagraphAfter()

# Paragraph F: hyperlink to techradar
$pE = $d.Paragraphs($targetIndex + 2)
$pE.Range.InsertParagraphAfter()
$pF = $d.Paragraphs($targetIndex + 3)
$rF = $pF.Range
$rF.Collapse(1)
$hF = $d.Hyperlinks.Add($rF, "http://www.techradar.com/news/top-5-best-disaster-recovery-services", "", "", "http://www.techradar.com/news/top-5-best-disaster-recovery-services")
$hF.Range.Font.Name = "Arial"
$hF.Range.Font.NameAscii = "Arial"
$hF.Range.Font.Size = 12

# Paragraph G: blank
$pF2 = $d.Paragraphs($targetIndex + 3)
$pF2.Range.InsertParagraphAfter()

# Paragraph H: hyperlink to zetta.net
$pG = $d.Paragraphs($targetIndex + 4)
$pG.Range.InsertParagraphAfter()
$pH = $d.Paragraphs($targetIndex + 5)
$rH = $pH.Range
$rH.Collapse(1)
$hH = $d.Hyperlinks.Add($rH, "https://www.zetta.net/zetta-disaster-recovery", "", "", "https://www.zetta.net/zetta-disaster-recovery")
$hH.Range.Font.Name = "Arial"
$hH.Range.Font.NameAscii = "Arial"
$hH.Range.Font.Size = 12

# Paragraph I: blank
$pH2 = $d.Paragraphs($targetIndex + 5)
$pH2.Range.InsertParagraphAfter()

# Paragraph J: hyperlink to planb.co.uk
$pI = $d.Paragraphs($targetIndex + 6)
$pI.Range.InsertParagraphAfter()
$pJ = $d.Paragraphs($targetIndex + 7)
$rJ = $pJ.Range
$rJ.Collapse(1)
$hJ = $d.Hyperlinks.Add($rJ, "https://www.planb.co.uk/services/", "", "", "https://www.planb.co.uk/services/")
$hJ.Range.Font.Name = "Arial"
$hJ.Range.Font.NameAscii = "Arial"
$hJ.Range.Font.Size = 12

Write-Output "done"
